$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.474.73'
$ws.Range('E2').Value = '  +3.79%  '
$ws.Range('D3').Value = '3.190.71'
$ws.Range('E3').Value = '  +4.89%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.66'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '634.87'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.234'
$ws.Range('E8').Value = '  +11.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.584'
$ws.Range('E9').Value = '  +5.62%  '
$ws.Range('D10').Value = '3.190.22'
$ws.Range('E10').Value = '  +4.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('E11').Value = '  +33.13%  '
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('E13').Value = '  +8.33%  '
$ws.Range('D14').Value = '3.780.38'
$ws.Range('E14').Value = '  +5.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000227'
$ws.Range('E15').Value = '  +16.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.85'
$ws.Range('E16').Value = '  +7.96%  '
$ws.Range('D17').Value = '79.266.34'
$ws.Range('E17').Value = '  +3.73%  '
$ws.Range('D18').Value = '3.192.65'
$ws.Range('E18').Value = '  +5.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.55'
$ws.Range('E19').Value = '  +8.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.10'
$ws.Range('E20').Value = '  +34.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.13'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '429.18'
$ws.Range('E22').Value = '  +13.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.03'
$ws.Range('E23').Value = '  +15.22%  '
$ws.Range('B24').Value = 'NEARProtocol'
$ws.Range('C24').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.81'
$ws.Range('E24').Value = '  +9.62%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.26'
$ws.Range('E25').Value = '  +12.88%  '
$ws.Range('D26').Value = '3.359.38'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '76.89'
$ws.Range('E27').Value = '  +4.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +5.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.05'
$ws.Range('E30').Value = '  +8.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  +4.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '525.53'
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('E35').Value = '  +26.96%  '
$ws.Range('E36').Value = '  +10.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +11.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.404'
$ws.Range('E39').Value = '  +4.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '164.70'
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '194.02'
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.03'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.49'
$ws.Range('E44').Value = '  +5.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.813'
$ws.Range('E45').Value = '  +11.32%  '
$ws.Range('E46').Value = '  +7.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.33'
$ws.Range('E47').Value = '  +4.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '42.98'
$ws.Range('E48').Value = '  +2.08%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.56'
$ws.Range('E49').Value = '  +4.71%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.80'
$ws.Range('E50').Value = '  +14.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.635'
$ws.Range('E51').Value = '  +4.70%  '
